$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for 2023Q4 (CVR numbers must stay text, like the existing
# column A entries, so quote-prefix them and then strip the transient
# formatting Excel attaches for the quote-prefix back off again).

# Row 40
$ws.Range("A40").Value = "'13893179"
$ws.Range("B40").Value = 2023
$ws.Range("C40").Value = 40172
$ws.Range("D40").Value = "Visma Løn"
$ws.Range("E40").NumberFormat = $ws.Range("E39").NumberFormat
$ws.Range("E40").Value = 45243
$ws.Range("G40").Value = "Zenegy"
$ws.Range("H40").Value = "2023Q4"
$ws.Range("I40").Value = "40000-60000"

# Row 41
$ws.Range("A41").Value = "'32762646"
$ws.Range("B41").Value = 2023
$ws.Range("C41").Value = 40892
$ws.Range("D41").Value = "EasyCruit"
$ws.Range("E41").NumberFormat = $ws.Range("E39").NumberFormat
$ws.Range("E41").Value = 45272
$ws.Range("H41").Value = "2023Q4"
$ws.Range("I41").Value = "40000-60000"

$ws.Range("A40:A41").ClearFormats()
